# Update "想去人数" (want-to-go count) values in both the "展览" sheet
# and the "全部类型" sheet, which both contain overlapping rows for the
# same events. Values only increase (refreshed scrape counts).

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 313
$wsExhibit.Range("F4").Value = 8231
$wsExhibit.Range("F5").Value = 6007
$wsExhibit.Range("F10").Value = 309
$wsExhibit.Range("F11").Value = 869
$wsExhibit.Range("F12").Value = 78

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 313
$wsAll.Range("F4").Value = 8231
$wsAll.Range("F5").Value = 6007
$wsAll.Range("F10").Value = 309
$wsAll.Range("F15").Value = 869
$wsAll.Range("F16").Value = 78
